$wb = $excel.ActiveWorkbook

# --- "About" sheet: update the "last updated" date in C1 (2/28/2024 -> 3/13/2024) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value2 = 45364

# --- "DPbES" sheet: set dispatch priority row 18 (municipal solid waste) from 0 to 1 across B:AE ---
$wsDP = $wb.Worksheets.Item("DPbES")
$wsDP.Range("B18:AE18").Value2 = 1

# --- Update the DPbES sheet view/selection state ---
$wsDP.Activate()
$wsDP.Range("E31").Select()
